$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): Right count 4 -> 5, Wrong marking -1 -> -1.2
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 (Total): Right total 56 -> 70, Wrong total 0 -> -0
$ws.Range("B12").Value = 70
$ws.Range("C12").Value = -0

# E12 label reflects new totals: "56/112" -> "70.0/140"
$ws.Range("E12").Value = "70.0/140"
